$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C2 (was "1234Z", now empty)
$ws.Range("C2").Value = ""

# G2 was text "150cc", now numeric 150
$ws.Range("G2").Value = 150

# G3 was text "1500cc", now numeric 1500
$ws.Range("G3").Value = 1500

# Update the active selection shown in the sheet view
$null = $ws.Range("D7").Select()
